$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "OSMO_DEF" header in column F, matching the style (bold, centered,
# bordered) already used by the existing B1:E1 headers.
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "OSMO_DEF"
$excel.CutCopyMode = 0
